$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row before row 51; this shifts old rows 51-164 down to 52-165
$ws.Rows("51").Insert()

# Populate the newly inserted row 51 with data (same categorical fields as the
# adjacent record, but a new date / volume / price observation)
$ws.Cells.Item(51, 1).Value = 4
$ws.Cells.Item(51, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(51, 3).Value = "Los Lagos"
$ws.Cells.Item(51, 4).Value = 44498
$ws.Cells.Item(51, 5).Value = 10
$ws.Cells.Item(51, 6).Value = "Fruta"
$ws.Cells.Item(51, 7).Value = 100102
$ws.Cells.Item(51, 8).Value = "Cítricos"
$ws.Cells.Item(51, 9).Value = 100102006
$ws.Cells.Item(51, 10).Value = "Pomelo"
$ws.Cells.Item(51, 11).Value = "Start Ruby"
$ws.Cells.Item(51, 12).Value = "Primera"
$ws.Cells.Item(51, 13).Value = 300
$ws.Cells.Item(51, 14).Value = 11000
$ws.Cells.Item(51, 15).Value = 12000
$ws.Cells.Item(51, 16).Value = 11500
$ws.Cells.Item(51, 17).Value = "`$/caja 14 kilos empedrada"
$ws.Cells.Item(51, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(51, 19).Value = 821
$ws.Cells.Item(51, 20).Value = 14

# Match the date format style used by the other date cells in column D
$ws.Cells.Item(51, 4).NumberFormat = $ws.Cells.Item(52, 4).NumberFormat
